$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 for Heap / Kth Largest Element in an Array (filled first)
$ws.Range("A7").Value = "Heap"
$ws.Range("B7").Value = "Kth Largest Element in an Array"
$ws.Range("C7").Value = "priority queue;"

# Match style of B7 to the "Sorting" category style (B5), since new Heap category uses the same accent style (s=3)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Rename existing category "D&C" -> "Divide & Conquer" on row 6 (edited last)
$ws.Range("A6").Value = "Divide & Conquer"

# Update selection to match final state
$ws.Range("B9").Select() | Out-Null
